$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column W (added for Waterland data) needs the same number format/style
# as the rest of the numeric data columns (B:V), which use format "0".
$ws.Range("W2:W8").NumberFormat = "0"

# Row 1: headers
$ws.Range("A1").Value = "year"
$ws.Range("B1").Value = "CRPland_acresk"
$ws.Range("C1").Value = "Cropland_acresk"
$ws.Range("D1").Value = "Federalland_acresk"
$ws.Range("E1").Value = "Forestland_acresk"
$ws.Range("F1").Value = "Pastureland_acresk"
$ws.Range("G1").Value = "Rangeland_acresk"
$ws.Range("H1").Value = "Ruralland_acresk"
$ws.Range("I1").Value = "Urbanland_acresk"
$ws.Range("J1").Value = "Waterland_acresk"
$ws.Range("K1").Value = "lccNA_acresk"
$ws.Range("L1").Value = "lccL1_acresk"
$ws.Range("M1").Value = "lccL2_acresk"
$ws.Range("N1").Value = "lccL3_acresk"
$ws.Range("O1").Value = "lccL4_acresk"
$ws.Range("P1").Value = "lccL5_acresk"
$ws.Range("Q1").Value = "lccL6_acresk"
$ws.Range("R1").Value = "lccL7_acresk"
$ws.Range("S1").Value = "lccL8_acresk"
$ws.Range("T1").Value = "lccL12_acresk"
$ws.Range("U1").Value = "lccL34_acresk"
$ws.Range("V1").Value = "lccL56_acresk"
$ws.Range("W1").Value = "lccL78_acresk"

# Row 2
$ws.Range("A2").Value = 1982
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 419711.70014646649
$ws.Range("D2").Value = 397162.801202178
$ws.Range("E2").Value = 408769.80048586428
$ws.Range("F2").Value = 130264.90003698319
$ws.Range("G2").Value = 417288.50006688386
$ws.Range("H2").Value = 62501.500497281551
$ws.Range("I2").Value = 50395.500235527754
$ws.Range("J2").Value = 49553.100219108164
$ws.Range("K2").Value = 518365.7019867152
$ws.Range("L2").Value = 29415.90003991127
$ws.Range("M2").Value = 285804.70027589798
$ws.Range("N2").Value = 287258.20012904704
$ws.Range("O2").Value = 201368.20010298491
$ws.Range("P2").Value = 34436.700000435114
$ws.Range("Q2").Value = 280848.30010822415
$ws.Range("R2").Value = 274476.90023375303
$ws.Range("S2").Value = 23673.200013324618
$ws.Range("T2").Value = 315220.60031580925
$ws.Range("U2").Value = 488626.40023203194
$ws.Range("V2").Value = 315285.00010865927
$ws.Range("W2").Value = 298150.10024707764

# Row 3
$ws.Range("A3").Value = 1987
$ws.Range("B3").Value = 13776.79998434335
$ws.Range("C3").Value = 405583.60013335943
$ws.Range("D3").Value = 397583.60121012479
$ws.Range("E3").Value = 410698.3004591614
$ws.Range("F3").Value = 126303.40002986044
$ws.Range("G3").Value = 411972.80007003248
$ws.Range("H3").Value = 62871.700501479208
$ws.Range("I3").Value = 56299.400283120573
$ws.Range("J3").Value = 50558.200218811631
$ws.Range("K3").Value = 525755.80204361677
$ws.Range("L3").Value = 29203.800039298832
$ws.Range("M3").Value = 283727.80025926977
$ws.Range("N3").Value = 285858.70011573285
$ws.Range("O3").Value = 200138.30009755492
$ws.Range("P3").Value = 34214.699999965727
$ws.Range("Q3").Value = 279714.00009713322
$ws.Range("R3").Value = 273535.6002253443
$ws.Range("S3").Value = 23499.100012376904
$ws.Range("T3").Value = 312931.60029856861
$ws.Range("U3").Value = 485997.00021328777
$ws.Range("V3").Value = 313928.70009709895
$ws.Range("W3").Value = 297034.7002377212

# Row 4
$ws.Range("A4").Value = 1992
$ws.Range("B4").Value = 34028.899985261261
$ws.Range("C4").Value = 381450.30009821802
$ws.Range("D4").Value = 399704.90122456104
$ws.Range("E4").Value = 410724.00043000281
$ws.Range("F4").Value = 124371.30002684146
$ws.Range("G4").Value = 408321.6000501439
$ws.Range("H4").Value = 63256.500508159399
$ws.Range("I4").Value = 63322.400341957808
$ws.Range("J4").Value = 50467.900225147605
$ws.Range("K4").Value = 534955.80212635547
$ws.Range("L4").Value = 28958.700038038194
$ws.Range("M4").Value = 281252.1002401337
$ws.Range("N4").Value = 283776.60009515285
$ws.Range("O4").Value = 198597.60008523613
$ws.Range("P4").Value = 33955.29999845475
$ws.Range("Q4").Value = 278276.50008951873
$ws.Range("R4").Value = 272506.50020512938
$ws.Range("S4").Value = 23368.700012274086
$ws.Range("T4").Value = 310210.8002781719
$ws.Range("U4").Value = 482374.20018038899
$ws.Range("V4").Value = 312231.80008797348
$ws.Range("W4").Value = 295875.20021740347

# Row 5
$ws.Range("A5").Value = 1997
$ws.Range("B5").Value = 32694.799986936152
$ws.Range("C5").Value = 375864.00007351488
$ws.Range("D5").Value = 400059.90122722834
$ws.Range("E5").Value = 411963.70039319992
$ws.Range("F5").Value = 119712.6999983117
$ws.Range("G5").Value = 406921.40004363656
$ws.Range("H5").Value = 63888.200500778854
$ws.Range("I5").Value = 73728.000436335802
$ws.Range("J5").Value = 50815.10023035109
$ws.Range("K5").Value = 546260.80223190039
$ws.Range("L5").Value = 28614.100037030876
$ws.Range("M5").Value = 278231.90020880103
$ws.Range("N5").Value = 281246.70006889105
$ws.Range("O5").Value = 196635.20006649196
$ws.Range("P5").Value = 33697.599996343255
$ws.Range("Q5").Value = 276691.10007810593
$ws.Range("R5").Value = 271132.90019249916
$ws.Range("S5").Value = 23137.500010229647
$ws.Range("T5").Value = 306846.00024583191
$ws.Range("U5").Value = 477881.90013538301
$ws.Range("V5").Value = 310388.70007444918
$ws.Range("W5").Value = 294270.40020272881

# Row 6
$ws.Range("A6").Value = 2002
$ws.Range("B6").Value = 31479.299977563322
$ws.Range("C6").Value = 367470.50004819036
$ws.Range("D6").Value = 401609.20124524087
$ws.Range("E6").Value = 412413.70034217089
$ws.Range("F6").Value = 118722.69999213517
$ws.Range("G6").Value = 406378.50002133101
$ws.Range("H6").Value = 64069.800497464836
$ws.Range("I6").Value = 82230.800529003143
$ws.Range("J6").Value = 51273.300237193704
$ws.Range("K6").Value = 557037.00235318393
$ws.Range("L6").Value = 28167.500035747886
$ws.Range("M6").Value = 275477.2001818344
$ws.Range("N6").Value = 278901.50004267693
$ws.Range("O6").Value = 195199.50004532933
$ws.Range("P6").Value = 33458.699993096292
$ws.Range("Q6").Value = 275395.90006704628
$ws.Range("R6").Value = 268967.10016188025
$ws.Range("S6").Value = 23043.400009498
$ws.Range("T6").Value = 303644.70021758229
$ws.Range("U6").Value = 474101.00008800626
$ws.Range("V6").Value = 308854.60006014258
$ws.Range("W6").Value = 292010.50017137825

# Row 7
$ws.Range("A7").Value = 2007
$ws.Range("B7").Value = 32578.799973286688
$ws.Range("C7").Value = 358786.00003223866
$ws.Range("D7").Value = 402130.50124900788
$ws.Range("E7").Value = 412162.40030286461
$ws.Range("F7").Value = 119433.79997787625
$ws.Range("G7").Value = 405568.40001321584
$ws.Range("H7").Value = 65602.700492627919
$ws.Range("I7").Value = 87770.900606565177
$ws.Range("J7").Value = 51614.300242610276
$ws.Range("K7").Value = 563588.90244240314
$ws.Range("L7").Value = 27857.100032843649
$ws.Range("M7").Value = 273569.40015505999
$ws.Range("N7").Value = 277538.1000245139
$ws.Range("O7").Value = 194130.30003011227
$ws.Range("P7").Value = 33300.2999914065
$ws.Range("Q7").Value = 274717.9000524655
$ws.Range("R7").Value = 267977.7001529038
$ws.Range("S7").Value = 22968.100008584559
$ws.Range("T7").Value = 301426.50018790364
$ws.Range("U7").Value = 471668.40005462617
$ws.Range("V7").Value = 308018.200043872
$ws.Range("W7").Value = 290945.80016148835

# Row 8
$ws.Range("A8").Value = 2012
$ws.Range("B8").Value = 23949.599979385734
$ws.Range("C8").Value = 361765.00001784414
$ws.Range("D8").Value = 402616.90125477314
$ws.Range("E8").Value = 412705.6002696529
$ws.Range("F8").Value = 121693.09997573495
$ws.Range("G8").Value = 404044.60000356287
$ws.Range("H8").Value = 66392.200491629541
$ws.Range("I8").Value = 90660.400652334094
$ws.Range("J8").Value = 51820.400245375931
$ws.Range("K8").Value = 567276.10249810666
$ws.Range("L8").Value = 27692.400029584765
$ws.Range("M8").Value = 272671.80014347285
$ws.Range("N8").Value = 276709.90000744909
$ws.Range("O8").Value = 193570.90002006292
$ws.Range("P8").Value = 33232.799990147352
$ws.Range("Q8").Value = 274218.80004697293
$ws.Range("R8").Value = 267305.10014606267
$ws.Range("S8").Value = 22970.000008434057
$ws.Range("T8").Value = 300364.20017305762
$ws.Range("U8").Value = 470280.80002751201
$ws.Range("V8").Value = 307451.60003712028
$ws.Range("W8").Value = 290275.10015449673
